$d = $word.ActiveDocument

# --- Change 1: merge "Socket para encaixe do " + "estojo " + "na " into a single run ---
$d.Content.Find.Execute("Socket para encaixe do estojo na ", $true, $false, $false, $false, $false, $true, 1, $false, "Socket para encaixe do estojo na ", 2) | Out-Null

# --- Change 2: merge "Cálculo da sensação térmica" + bookmark + ":" into a single run (drops the _GoBack bookmark) ---
$d.Content.Find.Execute("Cálculo da sensação térmica:", $true, $false, $false, $false, $false, $true, 1, $false, "Cálculo da sensação térmica:", 2) | Out-Null

# --- Change 3: fix "embulição" -> "ebulição" and re-split the runs, placing _GoBack at the edit point ---
$oldPara3 = "Elementos químicos – indicar o comportamento de diferentes elementos da tabela periódica de acordo com diferentes temperaturas. Por exemplo o ponto de fusão, ponto de embulição, etc."
$newPara3 = "Elementos químicos – indicar o comportamento de diferentes elementos da tabela periódica de acordo com diferentes temperaturas. Por exemplo o ponto de fusão, ponto de ebulição, etc."
$d.Content.Find.Execute($oldPara3, $true, $false, $false, $false, $false, $true, 1, $false, $newPara3, 2) | Out-Null

$rng = $d.Content.Duplicate
$rng.Find.Execute("Elementos químicos", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $rng.Start

# Split between "...Por exemp" and "lo o ponto..." (137 chars in) without leaving a stray bookmark.
$split1 = $d.Range($base + 137, $base + 137)
$tmp = $d.Bookmarks.Add("TempSplit", $split1)
$d.Bookmarks("TempSplit").Delete()

# Split between "...ponto de e" and "bulição, etc." (168 chars in), keeping the _GoBack bookmark there.
$split2 = $d.Range($base + 168, $base + 168)
$d.Bookmarks.Add("_GoBack", $split2) | Out-Null
